{"js": "// Replace the computed three-digit-by-one-digit multiplication answers\n// in the table cells with the new values from the commit.\nconst replacements = [\n  [\"411\u00d72=822\", \"541\u00d73=1623\"],\n  [\"577\u00d72=1154\", \"687\u00d72=1374\"],\n  [\"736\u00d73=2208\", \"704\u00d75=3520\"],\n  [\"345\u00d72=690\", \"640\u00d79=5760\"],\n  [\"391\u00d79=3519\", \"220\u00d72=440\"],\n  [\"834\u00d75=4170\", \"220\u00d78=1760\"],\n  [\"582\u00d79=5238\", \"707\u00d74=2828\"],\n  [\"114\u00d77=798\", \"583\u00d79=5247\"],\n  [\"669\u00d73=2007\", \"780\u00d74=3120\"],\n  [\"446\u00d76=2676\", \"339\u00d74=1356\"],\n  [\"711\u00d75=3555\", \"461\u00d75=2305\"],\n  [\"669\u00d76=4014\", \"410\u00d75=2050\"],\n  [\"571\u00d77=3997\", \"253\u00d77=1771\"],\n  [\"480\u00d77=3360\", \"578\u00d75=2890\"],\n  [\"278\u00d72=556\", \"168\u00d74=672\"],\n  [\"943\u00d76=5658\", \"688\u00d77=4816\"],\n  [\"777\u00d78=6216\", \"302\u00d78=2416\"],\n  [\"478\u00d74=1912\", \"840\u00d78=6720\"],\n  [\"556\u00d78=4448\", \"208\u00d76=1248\"],\n  [\"583\u00d73=1749\", \"687\u00d75=3435\"],\n  [\"210\u00d76=1260\", \"766\u00d79=6894\"],\n  [\"492\u00d73=1476\", \"675\u00d72=1350\"],\n  [\"446\u00d75=2230\", \"332\u00d78=2656\"],\n  [\"134\u00d79=1206\", \"269\u00d76=1614\"],\n  [\"663\u00d75=3315\", \"328\u00d73=984\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the computed three-digit-by-one-digit multiplication answers\n# in the table cells with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"411\u00d72=822\", \"541\u00d73=1623\"),\n    @(\"577\u00d72=1154\", \"687\u00d72=1374\"),\n    @(\"736\u00d73=2208\", \"704\u00d75=3520\"),\n    @(\"345\u00d72=690\", \"640\u00d79=5760\"),\n    @(\"391\u00d79=3519\", \"220\u00d72=440\"),\n    @(\"834\u00d75=4170\", \"220\u00d78=1760\"),\n    @(\"582\u00d79=5238\", \"707\u00d74=2828\"),\n    @(\"114\u00d77=798\", \"583\u00d79=5247\"),\n    @(\"669\u00d73=2007\", \"780\u00d74=3120\"),\n    @(\"446\u00d76=2676\", \"339\u00d74=1356\"),\n    @(\"711\u00d75=3555\", \"461\u00d75=2305\"),\n    @(\"669\u00d76=4014\", \"410\u00d75=2050\"),\n    @(\"571\u00d77=3997\", \"253\u00d77=1771\"),\n    @(\"480\u00d77=3360\", \"578\u00d75=2890\"),\n    @(\"278\u00d72=556\", \"168\u00d74=672\"),\n    @(\"943\u00d76=5658\", \"688\u00d77=4816\"),\n    @(\"777\u00d78=6216\", \"302\u00d78=2416\"),\n    @(\"478\u00d74=1912\", \"840\u00d78=6720\"),\n    @(\"556\u00d78=4448\", \"208\u00d76=1248\"),\n    @(\"583\u00d73=1749\", \"687\u00d75=3435\"),\n    @(\"210\u00d76=1260\", \"766\u00d79=6894\"),\n    @(\"492\u00d73=1476\", \"675\u00d72=1350\"),\n    @(\"446\u00d75=2230\", \"332\u00d78=2656\"),\n    @(\"134\u00d79=1206\", \"269\u00d76=1614\"),\n    @(\"663\u00d75=3315\", \"328\u00d73=984\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
